$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a full copy of the last existing data row (row 5) so the new
# row 6 inherits the same per-column styles (s="2" on VERSION, s="3" on the
# rest, no style on VIN/MAKE/MAKE_TEXT/MODEL_TEXT/BODYTYPE_TEXT/BODY_STYLE_CD).
$ws.Range("A5:AL5").Copy($ws.Range("A6:AL6"))

# Now overwrite the copied values with the new VIN record's data.
# (VIN in column A is written last so new shared strings are appended to
# the sharedStrings table in the same order as the reference edit.)
$ws.Range("B6").Value = "SYMBOL_2000_SS_TEST"
$ws.Range("C6").Value = 2017
$ws.Range("D6").Value = "Volkswagen"
$ws.Range("E6").Value = "Volkswagen"
$ws.Range("F6").Value = "Arteon"
$ws.Range("G6").Value = "Arteon SEL"
$ws.Range("H6").Value = 88888
$ws.Range("I6").Value = "WAG"
$ws.Range("J6").Value = "Coupe"
$ws.Range("K6").Value = "Sedan"
$ws.Range("L6").Value = "Coupe"
$ws.Range("M6").Value = "WAG"
$ws.Range("N6").Value = "8L V12"
$ws.Range("O6").Value = 12
$ws.Range("P6").Value = "G"
$ws.Range("Q6").Value = 214
$ws.Range("R6").Value = "4WD"
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = "000R"
$ws.Range("U6").Value = "DUAL AIR BAGS FRONT"
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = "4 WHEEL STANDARD"
$ws.Range("X6").Value = "STD"
$ws.Range("Y6").Value = "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"
$ws.Range("Z6").Value = "GC"
$ws.Range("AA6").Value = 35
$ws.Range("AB6").Value = 44
$ws.Range("AC6").Value = "S"
$ws.Range("AD6").Value = "Y"
$ws.Range("AE6").Value = "R"
$ws.Range("AF6").Value = "E"
$ws.Range("AG6").Value = "S"
$ws.Range("AH6").Value = "A"
$ws.Range("AI6").Value = 20000101
$ws.Range("AJ6").Value = "Y"
$ws.Range("AK6").Value = "Y"
$ws.Range("AL6").Value = "N"
$ws.Range("A6").Value = "7MSRP17H&V"

# Match the author's final selection state.
$ws.Range("B6").Select()
